# Adesão de verificação de pesquisa específica: nova coluna "DATA E HORÁRIO DA EXTRAÇÃO"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "DATA E HORÁRIO DA EXTRAÇÃO"
$ws.Columns.Item(8).AutoFit()

# mudanca de selecao ativa (prints/posicao no terminal ao rodar o script)
$ws.Range("D8").Select() | Out-Null

Write-Host "Coluna 'DATA E HORARIO DA EXTRACAO' adicionada com sucesso."
